# Add a new "PDP" worksheet right after "DataSet" (becomes the 2nd of 4
# sheets, pushing "Forms" and "Cards" one position later) and populate it
# with the PDP validation test data / helper rows, per:
#   "Hf and oxo pdp validation test case,test data, helper"

$wb = $excel.ActiveWorkbook

$dataSetSheet = $wb.Worksheets.Item("DataSet")

# Insert the new sheet immediately after "DataSet".
$pdp = $wb.Worksheets.Add($null, $dataSetSheet)
$pdp.Name = "PDP"

# ----- Header row (row 1) -----
$headers = @{
    "A1"  = "DataSet"
    "B1"  = "UserName"
    "C1"  = "Password"
    "D1"  = "Confirm Password"
    "E1"  = "FirstName"
    "F1"  = "LastName"
    "G1"  = "Email"
    "H1"  = "Color"
    "I1"  = "Colorproduct"
    "J1"  = "Products"
    "K1"  = "Street"
    "L1"  = "City"
    "M1"  = "Country"
    "N1"  = "State"
    "O1"  = "Region"
    "P1"  = "postcode"
    "Q1"  = "phone"
    "R1"  = "OTP Number"
    "S1"  = "Shippingmethods"
    "T1"  = "cardType"
    "U1"  = "cardNumber"
    "V1"  = "ExpMonthYear"
    "W1"  = "cvv"
    "X1"  = "Quantity"
    "Y1"  = "Discountcode"
    "Z1"  = "OXOAnswers"
    "AA1" = "Links"
    "AB1" = "productquantity"
}

foreach ($addr in $headers.Keys) {
    $cell = $pdp.Range($addr)
    $cell.Value = $headers[$addr]
    # Same yellow header fill used by the other sheets in the workbook.
    $cell.Interior.Color = 65535
}

# ----- Data row (row 2) -----
$pdp.Range("A2").Value = "Product"
$pdp.Range("H2").Value = "Tot Teal"
$pdp.Range("I2").Value = " On-the-Go Fork and Spoon Set"
$pdp.Range("AB2").Value = 3

# Match the selection/active cell recorded for the PDP sheet.
$pdp.Range("E5").Select()
